# Add a new "Hungary" market sheet to the FC Gallery Sounder panel workbook.
#
# Mirrors the real-world authoring flow: duplicate the last market sheet
# ("Slovakia") via "Move or Copy... > Create a copy", drop it in right after
# the source sheet, rename it, and fill in the two market-specific cells
# (Market name / NGC product code). The copy carries over all of Slovakia's
# formatting (column widths, merged cells, borders/fills) for free, so only
# the two changed cells need touching.

$wb = $excel.ActiveWorkbook
$slovakia = $wb.Worksheets.Item("Slovakia")

# Duplicate "Slovakia" and place the copy immediately after it (i.e. at the
# end of the tab strip).
$slovakia.Copy($null, $slovakia)
$hungary = $wb.Worksheets.Item($wb.Worksheets.Count)
$hungary.Name = "Hungary"

# Fill in the market-specific data.
$hungary.Range("B2").Value = "Hungary Market"
$hungary.Range("B4").Value = "NGC-4308/T3596/T3621"

# The product-code cell was typed in fresh (not copied with formatting), so
# it keeps the sheet's default style rather than Slovakia's bordered style.
$hungary.Range("B4").ClearFormats()

# Restore the UI selection state: Slovakia's selection becomes a full
# column/select-all selection (left over from copying it), while the new
# Hungary tab is the one left active, selected on its data entry cell.
$slovakia.Cells.Select()
$hungary.Activate()
$hungary.Range("B4").Select()
